# Generate Report for Handoff
#
# The localization pipeline finished handing the zh-cn content off for
# translation, so the status flips from "In Translation" to
# "Ready for handoff" and the corresponding handoff timestamps are
# refreshed on the Overview sheet and on each per-language detail sheet.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
# E2 = zh-cn status, F2 = de-de status, G2 = latest HO xliff generate date
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-28 14:38:48"

# --- zh-cn detail sheet ----------------------------------------------
# C2 = Status, H2 = Latest Handoff Datetime
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-28 14:38:43"

# --- de-de detail sheet ------------------------------------------------
# C2 = Status, H2 = Latest Handoff Datetime
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-28 14:38:48"

# --- Column widths -----------------------------------------------------
# "Ready for handoff" is longer than "In Translation", so the status
# columns widen to keep fitting the text (Overview!E:F and the Status
# column on each detail sheet).
$wsOverview.Columns.Item(5).ColumnWidth = 16.333333333333332
$wsOverview.Columns.Item(6).ColumnWidth = 16.333333333333332
$wsZhCn.Columns.Item(3).ColumnWidth = 16.333333333333332
$wsDeDe.Columns.Item(3).ColumnWidth = 16.333333333333332
